$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 13 (columns A:J) into new row 14, preserving exact value types.
# (Copy / PasteSpecial-values keeps text-typed, numeric-looking strings as text,
# unlike assigning .Value/.Value2 directly which would coerce them to
# numbers/booleans and lose the original shared-string typing.)
$ws.Range("A13:J13").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Update the view: scroll back so column A is the left-most visible column again
# and move the active selection to L17.
$ws.Range("A1").Select() | Out-Null
$ws.Range("L17").Select() | Out-Null
